$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up template/reference cells (scratch row 100) to carry the three column styles ---
$tmplLabel = $ws.Cells.Item(100, 1)
$tmplLabel.Value2 = "x"
$tmplLabel.Font.Bold = $true
$tmplLabel.Font.ColorIndex = -4105
$tmplLabel.WrapText = $false
$tmplLabel.VerticalAlignment = -4160

$tmplOrig = $ws.Cells.Item(100, 2)
$tmplOrig.Value2 = "x"
$tmplOrig.Font.Bold = $false
$tmplOrig.Font.ColorIndex = -4105
$tmplOrig.WrapText = $true
$tmplOrig.VerticalAlignment = -4160

$tmplMod = $ws.Cells.Item(100, 3)
$tmplMod.Value2 = "x"
$tmplMod.Font.Bold = $false
$tmplMod.Font.Color = 255
$tmplMod.WrapText = $true
$tmplMod.VerticalAlignment = -4160

function Set-Text($cell, $text, $tmpl) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $tmpl.Copy()
    $cell.PasteSpecial(-4122)
}

function Set-Label($r, $text) {
    Set-Text ($ws.Cells.Item($r, 1)) $text $tmplLabel
}
function Set-Orig($r, $text) {
    Set-Text ($ws.Cells.Item($r, 2)) $text $tmplOrig
}
function Set-Mod($r, $text) {
    Set-Text ($ws.Cells.Item($r, 3)) $text $tmplMod
}
function Clear-Row($r) {
    $ws.Rows.Item($r).ClearContents()
}

Clear-Row 1
Set-Orig 1 'Ementa atual:'
Set-Mod 1 'Ementa modificada (dados modificados em vermelho):'
$ws.Rows.Item(1).EntireRow.AutoFit()

Clear-Row 2
Set-Orig 2 'LOQ4053'
Set-Mod 2 'LOQ4053'
$ws.Rows.Item(2).EntireRow.AutoFit()

Clear-Row 3
Set-Label 3 'Nome:'
Set-Orig 3 ' Balanços de Massa e Energia'
Set-Mod 3 ' Balanços de Massa e Energia'
$ws.Rows.Item(3).EntireRow.AutoFit()

Clear-Row 4
Set-Label 4 'Name:'
Set-Orig 4 'Mass and Energy Balances'
Set-Mod 4 'Mass and Energy Balances'
$ws.Rows.Item(4).EntireRow.AutoFit()

Clear-Row 5
Set-Label 5 'Créditos-aula:'
Set-Orig 5 '2'
Set-Mod 5 '2'
$ws.Rows.Item(5).EntireRow.AutoFit()

Clear-Row 6
Set-Label 6 'Créditos-trabalho'
Set-Orig 6 '2'
Set-Mod 6 '2'
$ws.Rows.Item(6).EntireRow.AutoFit()

Clear-Row 7
Set-Label 7 'Carga horária:'
Set-Orig 7 '90 h'
Set-Mod 7 '90 h'
$ws.Rows.Item(7).EntireRow.AutoFit()

Clear-Row 8
Set-Label 8 'Ativação:'
Set-Orig 8 '01/01/2020'
Set-Mod 8 '01/01/2020'
$ws.Rows.Item(8).EntireRow.AutoFit()

Clear-Row 9
Set-Label 9 'Semestre ideal:'
Set-Orig 9 'EA-3,EB-3,EQD-2,EQN-3'
Set-Mod 9 'EA-3,EB-3,EQD-2,EQN-3'
$ws.Rows.Item(9).EntireRow.AutoFit()

Clear-Row 10
Set-Label 10 'Objetivos:'
Set-Orig 10 '5817045 - Elisângela de Jesus Cândido Moraes'
Set-Mod 10 '5817045 - Elisângela de Jesus Cândido Moraes'
$ws.Rows.Item(10).RowHeight = 60

Clear-Row 11
Set-Label 11 'Objectives:'
Set-Orig 11 'And introduce students to the basic setting to be used in all cases involving the conservation of mass and energy principles. This course provides the realization of global balances of mass and energy in different chemical processes highlighting the importance of using this methodology in the design and optimization of industrial chemical processes.'
Set-Mod 11 'And introduce students to the basic setting to be used in all cases involving the conservation of mass and energy principles. This course provides the realization of global balances of mass and energy in different chemical processes highlighting the importance of using this methodology in the design and optimization of industrial chemical processes.'
$ws.Rows.Item(11).RowHeight = 60

Clear-Row 12
Set-Label 12 'Docentes responsáveis:'
$ws.Rows.Item(12).EntireRow.AutoFit()

Clear-Row 13
Set-Label 13 'Programa resumido:'
Set-Orig 13 'Semestral'
Set-Mod 13 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

Clear-Row 14
Set-Label 14 'Short syllabus:'
Set-Orig 14 'Introduction to data calculations in Chemical Engineering; materials balance; Energy balances; materials and energy balances combined; Balances in processes in transient state.'
Set-Mod 14 'Introduction to data calculations in Chemical Engineering; materials balance; Energy balances; materials and energy balances combined; Balances in processes in transient state.'
$ws.Rows.Item(14).RowHeight = 60

Clear-Row 15
Set-Label 15 'Programa:'
Set-Orig 15 '01/01/2020'
Set-Mod 15 '01/01/2020'
$ws.Rows.Item(15).RowHeight = 120

Clear-Row 16
Set-Label 16 'Syllabus:'
Set-Orig 16 '1) Material Balancesa) Introduction to Material Balancesb) Material Balances that do not involve chemical reactions.c) Material Balances involving chemical reactions.d) Material Balances with recirculation (recycle and Bypass).2) Energy Balancesa) Definitions and concepts. Forms of energy, heat, enthalpy values of enthalpy and heat capacity.b) Balance of energy that do not involve chemical reactions.c) Energy balances involving chemical reactions.3) Mass and Balance of Power Combineda) Combined application of mass and energy balances in processes such as wetting, dissolving, mixing, etc.4) Mass Balance and Energy Processes in Transient Regime'
Set-Mod 16 '1) Material Balancesa) Introduction to Material Balancesb) Material Balances that do not involve chemical reactions.c) Material Balances involving chemical reactions.d) Material Balances with recirculation (recycle and Bypass).2) Energy Balancesa) Definitions and concepts. Forms of energy, heat, enthalpy values of enthalpy and heat capacity.b) Balance of energy that do not involve chemical reactions.c) Energy balances involving chemical reactions.3) Mass and Balance of Power Combineda) Combined application of mass and energy balances in processes such as wetting, dissolving, mixing, etc.4) Mass Balance and Energy Processes in Transient Regime'
$ws.Rows.Item(16).RowHeight = 120

Clear-Row 17
Set-Label 17 'Avaliação:'
$ws.Rows.Item(17).EntireRow.AutoFit()

Clear-Row 18
Set-Label 18 'Método:'
Set-Orig 18 '5817045 - Elisângela de Jesus Cândido Moraes'
Set-Mod 18 '5817045 - Elisângela de Jesus Cândido Moraes'
$ws.Rows.Item(18).RowHeight = 60

Clear-Row 19
Set-Label 19 'Critério:'
Set-Orig 19 'Provas escritas; -participação e conteúdo de trabalho e seminário;'
Set-Mod 19 'Provas escritas; -participação e conteúdo de trabalho e seminário;'
$ws.Rows.Item(19).RowHeight = 60

Clear-Row 20
Set-Label 20 'Norma de recuperação:'
Set-Orig 20 'Média Final = (Prova1 + 2xProva2 + Nota de Trabalho) / 4`nMédia final mínima de aprovação = 5,0'
Set-Mod 20 'Média Final = (Prova1 + 2xProva2 + Nota de Trabalho) / 4`nMédia final mínima de aprovação = 5,0'
$ws.Rows.Item(20).RowHeight = 60

Clear-Row 21
Set-Label 21 'Bibliografia:'
Set-Orig 21 '(Prova escrita + Média Final)/2`nNota Final mínima para aprovação= 5,0'
Set-Mod 21 '(Prova escrita + Média Final)/2`nNota Final mínima para aprovação= 5,0'
$ws.Rows.Item(21).RowHeight = 120

# Clean up scratch template row (do this BEFORE deleting row 22 so the
# row index of the template row is not shifted)
$ws.Rows.Item(100).EntireRow.Delete()

# Remove the now-superfluous row 22 (sheet now has only 21 rows of data)
$ws.Rows.Item(22).EntireRow.Delete()

Write-Host "done"